$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data added a new period column "Agosto.2021" (column BH) to the
# right of the existing last column "Mayo.2021" (column BG). For each data
# row the new column simply repeats the most recent known value (the one in
# column BG), matching how the other "latest period" columns were populated
# in previous updates of this series.

# Copy the full BG column (header + all 18 data rows) into BH so that both
# the values and the shared-string/number formatting are carried over.
$ws.Range("BG1:BG19").Copy()
$ws.Range("BH1:BH19").PasteSpecial(-4104)

# Also copy the header cell's formatting (bold, centered, bordered style)
# explicitly, since a values-only paste can leave the destination cell with
# default formatting.
$ws.Range("BG1").Copy()
$ws.Range("BH1").PasteSpecial(-4122)

# Set the new header label for the added period.
$ws.Range("BH1").Value = "Agosto.2021"

$excel.CutCopyMode = 0
